$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap whole rows for three country pairs that changed sort order ---
# (their own per-country data travels with the name)

function Swap-Rows($r1, $r2) {
    $cols = @("A","B","C","D","E","F","G","H")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Santa Lucia (188) / Belice (189) -> Belice now above Santa Lucia
Swap-Rows 188 189

# Namibia (194) / San Vicente y las Granadinas (195) -> San Vicente now above Namibia
Swap-Rows 194 195

# Burundi (198) / San Cristobal y Nieves (199) -> San Cristobal now above Burundi
Swap-Rows 198 199

# --- Update latest-day totals ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1187302
$ws.Range("C4").Value = 26528
$ws.Range("E4").Value = 940470
$ws.Range("G4").Value = 1125
$ws.Range("H4").Value = 68569

# Row 57: Argentina
$ws.Range("B57").Value = 4783
$ws.Range("C57").Value = 102
$ws.Range("E57").Value = 3183
$ws.Range("G57").Value = 9
$ws.Range("H57").Value = 246
